# Add "Wins" / "Losses" / "Ties" season-record columns (AD:AF) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-HeaderCell($addr, $text) {
    # Same look as the existing header cells (A1:AC1): bold font, thin box
    # border, centered horizontally, top-aligned vertically.
    $cell = $ws.Range($addr)
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous (thin box border)
}

# --- Header row (row 1) ---
Set-HeaderCell "AD1" "Wins"
Set-HeaderCell "AE1" "Losses"
Set-HeaderCell "AF1" "Ties"

# --- Data rows (2-71): season record is identical for every player row ---
$lastRow = 71
$ws.Range("AD2:AD$lastRow").Value = 71
$ws.Range("AE2:AE$lastRow").Value = 91
$ws.Range("AF2:AF$lastRow").Value = 0

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
